{"js": "// Replace each two-digit multiplication equation in the document with its\n// new value. Every \"old\" equation text is unique in the document, so a\n// simple ordered list of (old -> new) search/replace pairs is sufficient\n// and unambiguous.\nconst replacements = [\n  [\"71\u00d741=\", \"49\u00d720=\"],\n  [\"36\u00d774=\", \"80\u00d723=\"],\n  [\"74\u00d742=\", \"77\u00d740=\"],\n  [\"59\u00d779=\", \"93\u00d714=\"],\n  [\"91\u00d772=\", \"49\u00d756=\"],\n  [\"31\u00d770=\", \"87\u00d791=\"],\n  [\"96\u00d733=\", \"61\u00d725=\"],\n  [\"57\u00d792=\", \"20\u00d790=\"],\n  [\"44\u00d793=\", \"48\u00d741=\"],\n  [\"25\u00d713=\", \"62\u00d736=\"],\n  [\"65\u00d789=\", \"41\u00d737=\"],\n  [\"88\u00d747=\", \"71\u00d762=\"],\n  [\"86\u00d785=\", \"37\u00d717=\"],\n  [\"89\u00d717=\", \"67\u00d797=\"],\n  [\"44\u00d752=\", \"29\u00d722=\"],\n  [\"57\u00d775=\", \"83\u00d741=\"],\n  [\"39\u00d752=\", \"52\u00d770=\"],\n  [\"36\u00d748=\", \"13\u00d792=\"],\n  [\"19\u00d744=\", \"21\u00d763=\"],\n  [\"94\u00d786=\", \"58\u00d778=\"],\n  [\"99\u00d794=\", \"71\u00d733=\"],\n  [\"92\u00d740=\", \"74\u00d751=\"],\n  [\"56\u00d723=\", \"64\u00d784=\"],\n  [\"24\u00d790=\", \"56\u00d722=\"],\n  [\"33\u00d711=\", \"77\u00d760=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication equation in the document with its\n# new value. Every \"old\" equation text is unique in the document, so a\n# simple ordered list of (old -> new) find/replace pairs is sufficient and\n# unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"71\u00d741=\", \"49\u00d720=\"),\n  @(\"36\u00d774=\", \"80\u00d723=\"),\n  @(\"74\u00d742=\", \"77\u00d740=\"),\n  @(\"59\u00d779=\", \"93\u00d714=\"),\n  @(\"91\u00d772=\", \"49\u00d756=\"),\n  @(\"31\u00d770=\", \"87\u00d791=\"),\n  @(\"96\u00d733=\", \"61\u00d725=\"),\n  @(\"57\u00d792=\", \"20\u00d790=\"),\n  @(\"44\u00d793=\", \"48\u00d741=\"),\n  @(\"25\u00d713=\", \"62\u00d736=\"),\n  @(\"65\u00d789=\", \"41\u00d737=\"),\n  @(\"88\u00d747=\", \"71\u00d762=\"),\n  @(\"86\u00d785=\", \"37\u00d717=\"),\n  @(\"89\u00d717=\", \"67\u00d797=\"),\n  @(\"44\u00d752=\", \"29\u00d722=\"),\n  @(\"57\u00d775=\", \"83\u00d741=\"),\n  @(\"39\u00d752=\", \"52\u00d770=\"),\n  @(\"36\u00d748=\", \"13\u00d792=\"),\n  @(\"19\u00d744=\", \"21\u00d763=\"),\n  @(\"94\u00d786=\", \"58\u00d778=\"),\n  @(\"99\u00d794=\", \"71\u00d733=\"),\n  @(\"92\u00d740=\", \"74\u00d751=\"),\n  @(\"56\u00d723=\", \"64\u00d784=\"),\n  @(\"24\u00d790=\", \"56\u00d722=\"),\n  @(\"33\u00d711=\", \"77\u00d760=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  # wdFindContinue = 1, wdReplaceAll = 2\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
